# Apply updated crypto price/volume data per the source diff.
# Note: several "Price" values look numeric (e.g. "1.007") but are stored as
# plain text in the workbook. Prefixing with an apostrophe forces Excel to
# keep them as text instead of auto-converting them to numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"

$ws.Range("D2").Value = '26.286.36'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '1.676.76'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = $apos + '217.54'
$ws.Range("D6").Value = $apos + '0.5295'
$ws.Range("E6").Value = '  +4.38%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = $apos + '0.2688'
$ws.Range("E8").Value = '  +2.23%  '
$ws.Range("D9").Value = $apos + '0.06471'
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("D10").Value = $apos + '21.91'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").Value = $apos + '0.07520'
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = $apos + '4.511'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.662.12'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").Value = $apos + '0.5776'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = $apos + '0.000008491'
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = $apos + '64.65'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '26.320.22'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = $apos + '4.921'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = $apos + '1.007'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = $apos + '10.86'
$ws.Range("E20").Value = '  +1.98%  '
$ws.Range("D21").Value = $apos + '189.99'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").Value = $apos + '6.193'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = $apos + '1.008'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = $apos + '144.70'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = $apos + '0.1274'
$ws.Range("E25").Value = '  +7.43%  '
$ws.Range("D26").Value = $apos + '7.804'
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("D28").Value = $apos + '0.06482'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").Value = $apos + '1.362'
$ws.Range("E29").Value = '  +4.19%  '
$ws.Range("D30").Value = $apos + '1.319'
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").Value = $apos + '3.586'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("D32").Value = $apos + '3.585'
$ws.Range("E32").Value = '  +2.69%  '
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").Value = $apos + '1.029'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = $apos + '0.6188'
$ws.Range("E35").Value = '  +2.36%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").Value = $apos + '2.735'
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").Value = $apos + '6.281'
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("D39").Value = '1.116.07'
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("D40").Value = $apos + '0.01620'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("D41").Value = $apos + '0.8705'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").Value = $apos + '100.51'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '1.826.95'
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("D45").Value = $apos + '0.00000000108'
$ws.Range("E45").Value = '  -5.12%  '
$ws.Range("D46").Value = $apos + '56.93'
$ws.Range("E46").Value = '  +1.66%  '
$ws.Range("D47").Value = $apos + '1.009'
$ws.Range("E47").Value = '  +0.65%  '
$ws.Range("D48").Value = $apos + '8.158'
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").Value = $apos + '0.05264'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").Value = $apos + '6.047'
$ws.Range("E51").Value = '  +2.02%  '
